$d = $word.ActiveDocument

$d.Content.Find.Execute("Categories and tags", $false, $false, $false, $false, $false, $true, 1, $false, "Standards", 2) | Out-Null
